$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.818.43"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "1.856.13"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.24"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5058"
$ws.Range("E7").Value = "  -0.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3639"
$ws.Range("E8").Value = "  -3.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07159"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8903"
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.62"
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "1.857.99"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07442"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.43"
$ws.Range("E14").Value = "  +3.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.225"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008488"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.02"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "26.866.79"
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("D22").Value = "2.093.13"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  -3.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.418"
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.30"
$ws.Range("E25").Value = "  -2.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.798"
$ws.Range("E26").Value = "  -2.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.82"
$ws.Range("E28").Value = "  -3.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.07"
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.629"
$ws.Range("E30").Value = "  -2.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.653"
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09220"
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05080"
$ws.Range("E33").Value = "  -1.65%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.976"
$ws.Range("E34").Value = "  -3.79%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7435"
$ws.Range("E35").Value = "  -1.48%  "
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.252"
$ws.Range("E37").Value = "  +7.29%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.101"
$ws.Range("E38").Value = "  +2.05%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.502"
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01986"
$ws.Range("E40").Value = "  -2.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5317"
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "118.51"
$ws.Range("E42").Value = "  +3.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.459"
$ws.Range("E43").Value = "  -2.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.372"
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1457"
$ws.Range("E45").Value = "  -1.77%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4634"
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.986"
$ws.Range("E48").Value = "  -1.26%  "
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.77"
$ws.Range("E51").Value = "  -3.74%  "
